# "remove unused define for item"
#
# The ItemSubType column (C) held a now-unused sub-type "define" for most
# rows; clear it to 0. For the HOLY_WATER_* potion rows (B46:C61) the
# ItemType/ItemSubType pair is renumbered from (5,0) to (2,7).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 13-45: ItemSubType (column C) -> 0
for ($r = 13; $r -le 45; $r++) {
    $ws.Cells.Item($r, 3).Value = 0
}

# Rows 46-61: ItemType (column B) -> 2, ItemSubType (column C) -> 7
for ($r = 46; $r -le 61; $r++) {
    $ws.Cells.Item($r, 2).Value = 2
    $ws.Cells.Item($r, 3).Value = 7
}

# Rows 132-140: ItemSubType (column C) -> 0
for ($r = 132; $r -le 140; $r++) {
    $ws.Cells.Item($r, 3).Value = 0
}

$ws.Range("D131").Select()
